$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert 3 new blank rows before row 33 so the detail table grows
#     from 18 rows (16-33) to 21 rows (16-36). This correctly shifts the
#     closing/footer rows (old 38-39) down to (new 41-42) and updates the
#     sheet dimension + merged cell references automatically.
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# --- Step 2: Stamp the correct cell style onto the 3 newly inserted rows by
#     copying the formatting (and content, to be overwritten below) from the
#     last "normal" detail row (row 32, style ids 15/16/17/16/18/18/19/19/20).
$ws.Range("B32:J32").Copy($ws.Range("B33:J33"))
$ws.Range("B32:J32").Copy($ws.Range("B34:J34"))
$ws.Range("B32:J32").Copy($ws.Range("B35:J35"))

# --- Step 3: Rewrite the whole worker/period detail table (rows 16-36) with
#     the new data: reorganized by worker (instead of by period), including
#     the new period 2507 and recalculated "Valor Mora" amounts.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143401756"
$ws.Range("D16").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 119467
$ws.Range("G16").Value = 3200000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143401756"
$ws.Range("D17").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 128000
$ws.Range("G17").Value = 3200000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143401756"
$ws.Range("D18").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 128000
$ws.Range("G18").Value = 3200000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143401756"
$ws.Range("D19").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 128000
$ws.Range("G19").Value = 3200000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143401756"
$ws.Range("D20").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 128000
$ws.Range("G20").Value = 3200000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143401756"
$ws.Range("D21").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E21").Value = "2502"
$ws.Range("F21").Value = 128000
$ws.Range("G21").Value = 3200000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143401756"
$ws.Range("D22").Value = "LUIS ANGEL LUNA ESCORCIA"
$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 128000
$ws.Range("G22").Value = 3200000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1007230501"
$ws.Range("D23").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 52267
$ws.Range("G23").Value = 1900000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1007230501"
$ws.Range("D24").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E24").Value = "2506"
$ws.Range("F24").Value = 76000
$ws.Range("G24").Value = 1900000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1007230501"
$ws.Range("D25").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E25").Value = "2505"
$ws.Range("F25").Value = 76000
$ws.Range("G25").Value = 1900000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1007230501"
$ws.Range("D26").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E26").Value = "2504"
$ws.Range("F26").Value = 76000
$ws.Range("G26").Value = 1900000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1007230501"
$ws.Range("D27").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E27").Value = "2503"
$ws.Range("F27").Value = 76000
$ws.Range("G27").Value = 1900000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1007230501"
$ws.Range("D28").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E28").Value = "2502"
$ws.Range("F28").Value = 76000
$ws.Range("G28").Value = 1900000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1007230501"
$ws.Range("D29").Value = "SEBASTIAN DE JESUS RUIZ AVILA"
$ws.Range("E29").Value = "2501"
$ws.Range("F29").Value = 76000
$ws.Range("G29").Value = 1900000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1143351761"
$ws.Range("D30").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E30").Value = "2507"
$ws.Range("F30").Value = 74667
$ws.Range("G30").Value = 2000000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1143351761"
$ws.Range("D31").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E31").Value = "2506"
$ws.Range("F31").Value = 80000
$ws.Range("G31").Value = 2000000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1143351761"
$ws.Range("D32").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E32").Value = "2505"
$ws.Range("F32").Value = 80000
$ws.Range("G32").Value = 2000000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1143351761"
$ws.Range("D33").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E33").Value = "2504"
$ws.Range("F33").Value = 80000
$ws.Range("G33").Value = 2000000

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1143351761"
$ws.Range("D34").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E34").Value = "2503"
$ws.Range("F34").Value = 80000
$ws.Range("G34").Value = 2000000

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1143351761"
$ws.Range("D35").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E35").Value = "2502"
$ws.Range("F35").Value = 80000
$ws.Range("G35").Value = 2000000

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1143351761"
$ws.Range("D36").Value = "JAIRO JESUS CABARCAS MARTINEZ"
$ws.Range("E36").Value = "2501"
$ws.Range("F36").Value = 80000
$ws.Range("G36").Value = 2000000

# --- Step 4: Update the summary figures above the table.
$ws.Range("E11").Value = 1950401   # Valor Mora (sum of all F16:F36)
$ws.Range("F13").Value = 7         # Cant. Periodos (was 6, now 7: 2501-2507)
